# Refresh "outputs-HGR-r202" prediction export: replace the placeholder
# genome-length values in column B with the freshly computed lengths, and
# re-point column A at the MAG list from the current (post-filter) run —
# three genomes dropped out of this rerun, so every row from the first
# drop onward shifts up to the next surviving name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "even_MAG-GUT10153.fa", "even_MAG-GUT103.fa", "even_MAG-GUT10562.fa", "even_MAG-GUT10734.fa", "even_MAG-GUT11041.fa", "even_MAG-GUT11308.fa",
    "even_MAG-GUT11426.fa", "even_MAG-GUT11521.fa", "even_MAG-GUT11638.fa", "even_MAG-GUT11820.fa", "even_MAG-GUT11829.fa", "even_MAG-GUT11847.fa",
    "even_MAG-GUT11972.fa", "even_MAG-GUT11977.fa", "even_MAG-GUT12030.fa", "even_MAG-GUT12049.fa", "even_MAG-GUT12051.fa", "even_MAG-GUT12063.fa",
    "even_MAG-GUT12082.fa", "even_MAG-GUT12095.fa", "even_MAG-GUT12230.fa", "even_MAG-GUT12257.fa", "even_MAG-GUT12300.fa", "even_MAG-GUT12797.fa",
    "even_MAG-GUT13150.fa", "even_MAG-GUT13451.fa", "even_MAG-GUT13856.fa", "even_MAG-GUT14027.fa", "even_MAG-GUT14364.fa", "even_MAG-GUT14645.fa",
    "even_MAG-GUT14705.fa", "even_MAG-GUT14711.fa", "even_MAG-GUT1694.fa", "even_MAG-GUT1761.fa", "even_MAG-GUT18031.fa", "even_MAG-GUT18037.fa",
    "even_MAG-GUT19155.fa", "even_MAG-GUT19423.fa", "even_MAG-GUT20126.fa", "even_MAG-GUT20293.fa", "even_MAG-GUT20461.fa", "even_MAG-GUT22496.fa",
    "even_MAG-GUT22619.fa", "even_MAG-GUT22830.fa", "even_MAG-GUT24184.fa", "even_MAG-GUT24574.fa", "even_MAG-GUT24670.fa", "even_MAG-GUT27127.fa",
    "even_MAG-GUT27129.fa", "even_MAG-GUT27169.fa", "even_MAG-GUT27309.fa", "even_MAG-GUT27329.fa", "even_MAG-GUT27936.fa", "even_MAG-GUT27968.fa",
    "even_MAG-GUT2841.fa", "even_MAG-GUT2847.fa", "even_MAG-GUT2867.fa", "even_MAG-GUT2873.fa", "even_MAG-GUT28915.fa", "even_MAG-GUT31002.fa",
    "even_MAG-GUT31343.fa", "even_MAG-GUT31540.fa", "even_MAG-GUT317.fa", "even_MAG-GUT32169.fa", "even_MAG-GUT32170.fa", "even_MAG-GUT32343.fa",
    "even_MAG-GUT32351.fa", "even_MAG-GUT32588.fa", "even_MAG-GUT32929.fa", "even_MAG-GUT33270.fa", "even_MAG-GUT33329.fa", "even_MAG-GUT35732.fa",
    "even_MAG-GUT35747.fa", "even_MAG-GUT35871.fa", "even_MAG-GUT35943.fa", "even_MAG-GUT36571.fa", "even_MAG-GUT36733.fa", "even_MAG-GUT36772.fa",
    "even_MAG-GUT36796.fa", "even_MAG-GUT36799.fa", "even_MAG-GUT36814.fa", "even_MAG-GUT37401.fa", "even_MAG-GUT379.fa", "even_MAG-GUT37917.fa",
    "even_MAG-GUT38735.fa", "even_MAG-GUT39139.fa", "even_MAG-GUT40221.fa", "even_MAG-GUT41097.fa", "even_MAG-GUT41378.fa", "even_MAG-GUT41860.fa",
    "even_MAG-GUT42780.fa", "even_MAG-GUT42852.fa", "even_MAG-GUT42865.fa", "even_MAG-GUT43110.fa", "even_MAG-GUT43251.fa", "even_MAG-GUT43377.fa",
    "even_MAG-GUT43378.fa", "even_MAG-GUT43648.fa", "even_MAG-GUT43957.fa", "even_MAG-GUT44774.fa", "even_MAG-GUT44944.fa", "even_MAG-GUT45122.fa",
    "even_MAG-GUT45214.fa", "even_MAG-GUT45263.fa", "even_MAG-GUT45396.fa", "even_MAG-GUT4552.fa", "even_MAG-GUT45670.fa", "even_MAG-GUT4585.fa",
    "even_MAG-GUT4600.fa", "even_MAG-GUT46167.fa", "even_MAG-GUT4634.fa", "even_MAG-GUT46378.fa", "even_MAG-GUT4651.fa", "even_MAG-GUT47205.fa",
    "even_MAG-GUT47800.fa", "even_MAG-GUT48077.fa", "even_MAG-GUT48276.fa", "even_MAG-GUT4969.fa", "even_MAG-GUT4979.fa", "even_MAG-GUT50478.fa",
    "even_MAG-GUT50483.fa", "even_MAG-GUT50508.fa", "even_MAG-GUT51557.fa", "even_MAG-GUT51559.fa", "even_MAG-GUT5249.fa", "even_MAG-GUT52992.fa",
    "even_MAG-GUT53689.fa", "even_MAG-GUT5375.fa", "even_MAG-GUT54574.fa", "even_MAG-GUT54831.fa", "even_MAG-GUT55210.fa", "even_MAG-GUT55743.fa",
    "even_MAG-GUT56345.fa", "even_MAG-GUT56417.fa", "even_MAG-GUT57158.fa", "even_MAG-GUT57729.fa", "even_MAG-GUT5848.fa", "even_MAG-GUT58938.fa",
    "even_MAG-GUT59149.fa", "even_MAG-GUT593.fa", "even_MAG-GUT59579.fa", "even_MAG-GUT59599.fa", "even_MAG-GUT62658.fa", "even_MAG-GUT63164.fa",
    "even_MAG-GUT63198.fa", "even_MAG-GUT63214.fa", "even_MAG-GUT63219.fa", "even_MAG-GUT63286.fa", "even_MAG-GUT63373.fa", "even_MAG-GUT63496.fa",
    "even_MAG-GUT63586.fa", "even_MAG-GUT63602.fa", "even_MAG-GUT6384.fa", "even_MAG-GUT6445.fa", "even_MAG-GUT6489.fa", "even_MAG-GUT65795.fa",
    "even_MAG-GUT66097.fa", "even_MAG-GUT66701.fa", "even_MAG-GUT67224.fa", "even_MAG-GUT68072.fa", "even_MAG-GUT68311.fa", "even_MAG-GUT68785.fa",
    "even_MAG-GUT6955.fa", "even_MAG-GUT6968.fa", "even_MAG-GUT7012.fa", "even_MAG-GUT70404.fa", "even_MAG-GUT70913.fa", "even_MAG-GUT71577.fa",
    "even_MAG-GUT722.fa", "even_MAG-GUT72293.fa", "even_MAG-GUT72929.fa", "even_MAG-GUT73862.fa", "even_MAG-GUT74916.fa", "even_MAG-GUT76034.fa",
    "even_MAG-GUT761.fa", "even_MAG-GUT76518.fa", "even_MAG-GUT77597.fa", "even_MAG-GUT77633.fa", "even_MAG-GUT78579.fa", "even_MAG-GUT78910.fa",
    "even_MAG-GUT78923.fa", "even_MAG-GUT80330.fa", "even_MAG-GUT81123.fa", "even_MAG-GUT81409.fa", "even_MAG-GUT81646.fa", "even_MAG-GUT81936.fa",
    "even_MAG-GUT82089.fa", "even_MAG-GUT82177.fa", "even_MAG-GUT82314.fa", "even_MAG-GUT82472.fa", "even_MAG-GUT82505.fa", "even_MAG-GUT82571.fa",
    "even_MAG-GUT8267.fa", "even_MAG-GUT82998.fa", "even_MAG-GUT83643.fa", "even_MAG-GUT838.fa", "even_MAG-GUT84166.fa", "even_MAG-GUT84304.fa",
    "even_MAG-GUT85906.fa", "even_MAG-GUT86112.fa", "even_MAG-GUT86504.fa"
)

$values = @(
    940.71008313601578, 1020.107469590683, 512.59965668498512, 1567.4078853265389, 994.19542152882877, 817.81887403784185,
    894.5023493622557, 306.44659084786787, 235.05012499429125, 413.51132709682548, 751.19478832031245, 281.15848284538021,
    762.98239432514947, 376.50505377510132, 1480.1278573629938, 709.92825973824688, 812.91908019221296, 414.6520297883132,
    703.25049747096682, 409.89257690578989, 554.46798947148409, 791.24581057149339, 1536.7456354601679, 1040.0651044450926,
    943.75048948242693, 847.75242510088765, 727.85255884362846, 1438.845430706572, 1040.684799605526, 510.56343743872873,
    703.21353943551549, 528.40728598811484, 943.99178712000048, 1159.6228509874179, 683.11920172131613, 804.42185195584671,
    801.96843022655776, 975.8609936856376, 966.23546594888762, 1273.2742951763253, 986.35925123418451, 1065.1732582564684,
    960.52985855154316, 983.32316199007141, 1100.9063996944756, 1097.8362941436485, 1162.7189766254087, 1299.4942711528142,
    1528.7374166993548, 1247.8104549188472, 932.50266566448454, 1074.4303254708987, 964.04374513434459, 838.71130262380393,
    472.10560982204976, 884.52779146576245, 674.77323496908934, 953.24186982032234, 398.68740206826993, 822.5922229406201,
    1146.5253018781366, 941.63658990784779, 570.17120413063321, 204.73591427604879, 394.05295232651014, 1626.1966164675587,
    1666.5079390925569, 676.60501438768006, 824.93861879722681, 524.46765413286835, 605.22922377272903, 1061.7847420968126,
    628.34972078436772, 687.3740865151276, 911.51755415316939, 902.42933592127793, 1570.7264501421298, 873.81703830605989,
    991.140195591478, 1014.487903238396, 436.44556659896494, 653.28021878638742, 701.64888996910054, 1400.5953972205471,
    981.0797891882346, 659.72789067990038, 341.79052099911684, 805.32050190498035, 1394.7319550367931, 957.42653041126437,
    1029.0247797799941, 1329.4615182883067, 1775.9210126133744, 929.98623793345064, 762.98201961308064, 1401.1603609688061,
    1548.9814336897632, 864.19178926949144, 960.24102853197428, 767.66786742605223, 1073.3919253257222, 978.27913903164699,
    1101.1346474834186, 662.4679227243164, 870.77158069772463, 1305.5301167520772, 646.58275180579471, 1050.6253443937333,
    1521.54002131229, 1039.8572284896486, 1175.9882188885267, 718.64877928389137, 946.54402234876193, 603.79184191501918,
    1049.7001492335689, 1417.6939229598327, 884.06627231085758, 998.60665433740655, 899.09793383453143, 893.86176051540428,
    1652.9846716485993, 1200.4657423024014, 991.0835857966556, 956.52597642629826, 705.22378743076683, 965.52541704113287,
    1322.7959454895363, 1011.0236942060831, 1277.5773442770405, 1264.2288413608255, 757.20939447556634, 1340.0449194905152,
    1036.3627255987933, 1041.1205803310218, 705.13432795678045, 1334.3458243548471, 1209.1740088957097, 1107.8674152404406,
    1505.3042380701456, 995.27344335345174, 541.88093790056791, 796.14689923506216, 1017.6493709106624, 644.45936542207573,
    240.44300053812918, 229.79213156914597, 160.36366628164319, 202.57708733591068, 187.8272754387994, 210.0657832764681,
    213.53815116978188, 206.13603754838999, 665.79336910060033, 669.10009581143936, 819.45598948323573, 341.79052099911684,
    1393.4031261667453, 626.70490282627065, 1223.0679213573126, 650.90027976206898, 843.64830687173207, 1255.0678160484233,
    585.92172127489039, 484.39456939896428, 734.89107627895305, 1192.3812251245572, 282.82215759820662, 588.05285681662303,
    819.55090240428649, 1025.3617287269969, 696.36479892337343, 817.83371921907656, 1542.3300636430756, 816.78791199640659,
    909.27144631730312, 1005.3535197822789, 816.59364680809995, 1124.2575606191519, 1382.8759034850746, 678.8979574815105,
    621.84749897212987, 978.28326742838613, 819.94061780537663, 903.2933501435682, 853.80684288159, 717.71780466036012,
    573.66971693705887, 903.53176372221651, 589.55259047158984, 745.52318155683884, 769.76263111090589, 879.29100461188341,
    651.98994849802034, 962.84450733129324, 626.54148063508228, 932.37971746405697, 587.18530643348663, 776.16463472056921,
    843.21434942071892, 816.72858786416759, 1292.2860709773247
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $names[$i]
    $ws.Cells.Item($row, 2).Value2 = $values[$i]
}
